$wb = $excel.ActiveWorkbook

# Rename original sheet to "SerDes"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "SerDes"

# Add a new worksheet after "SerDes" and name it "Sheet1"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

$ws2.Range("A3").Value = "BS_ADC_PwrDn"
$ws2.Range("B3").Value = 11
$ws2.Range("C3").Value = "DdsIO_Reset  "
$ws2.Range("K3").Value = "            "
$ws2.Range("A4").Value = "BS_RefPll]   "
$ws2.Range("B4").Value = 15
$ws2.Range("C4").Value = "DdsTxEn      "
$ws2.Range("A5").Value = "BS_PllOut]   "
$ws2.Range("B5").Value = 17
$ws2.Range("C5").Value = "DdsSyncClk   "
$ws2.Range("A6").Value = "BS_I2C_SCl]  "
$ws2.Range("B6").Value = 19
$ws2.Range("C6").Value = "DdsIO_Update "
$ws2.Range("A7").Value = "BS_I2C_SDa]  "
$ws2.Range("B7").Value = 21
$ws2.Range("C7").Value = "DdsReset     "
$ws2.Range("A8").Value = "BS_DAC_Sel_n]"
$ws2.Range("B8").Value = 23
$ws2.Range("C8").Value = "DdsCS_n      "
$ws2.Range("A9").Value = "BS_DAC_SClk] "
$ws2.Range("B9").Value = 25
$ws2.Range("C9").Value = "DdsSClk      "
$ws2.Range("A10").Value = "BS_DAC_MOSI] "
$ws2.Range("B10").Value = 27
$ws2.Range("C10").Value = "DdsMosi      "
$ws2.Range("A11").Value = "BS_ADC_SE  "
$ws2.Range("B11").Value = 29
$ws2.Range("C11").Value = "DdsMiso      "
$ws2.Range("A13").Value = "BS_Ovf"
$ws2.Range("B13").Value = 2
$ws2.Range("C13").Value = "IF/BS_n"
$ws2.Range("A14").Value = "BS_D13"
$ws2.Range("B14").Value = 4
$ws2.Range("C14").Value = "DdsData17"
$ws2.Range("A15").Value = "BS_D12"
$ws2.Range("B15").Value = 6
$ws2.Range("C15").Value = "DdsData16"
$ws2.Range("A16").Value = "BS_D11"
$ws2.Range("B16").Value = 8
$ws2.Range("C16").Value = "DdsData15"
$ws2.Range("A17").Value = "BS_D10"
$ws2.Range("B17").Value = 10
$ws2.Range("C17").Value = "DdsData14"
$ws2.Range("A18").Value = "BS_D9"
$ws2.Range("B18").Value = 12
$ws2.Range("C18").Value = "DdsData13"
$ws2.Range("A19").Value = "BS_D8"
$ws2.Range("B19").Value = 14
$ws2.Range("C19").Value = "DdsData12"
$ws2.Range("F19").Value = "BS_Ovf]      "
$ws2.Range("A20").Value = "BS_D7"
$ws2.Range("B20").Value = 16
$ws2.Range("C20").Value = "DdsData11"
$ws2.Range("A21").Value = "BS_D6"
$ws2.Range("B21").Value = 18
$ws2.Range("C21").Value = "DdsData10"
$ws2.Range("A22").Value = "BS_D5"
$ws2.Range("B22").Value = 20
$ws2.Range("C22").Value = "DdsData9 "
$ws2.Range("A23").Value = "BS_D4"
$ws2.Range("B23").Value = 22
$ws2.Range("C23").Value = "DdsData8 "
$ws2.Range("A24").Value = "BS_D3"
$ws2.Range("B24").Value = 24
$ws2.Range("C24").Value = "DdsData7 "
$ws2.Range("A25").Value = "BS_D2"
$ws2.Range("B25").Value = 26
$ws2.Range("C25").Value = "DdsData6 "
$ws2.Range("A26").Value = "BS_D1"
$ws2.Range("B26").Value = 28
$ws2.Range("C26").Value = "DdsData5 "
$ws2.Range("A27").Value = "BS_D0"
$ws2.Range("B27").Value = 30
$ws2.Range("C27").Value = "DdsData4 "
$ws2.Range("B28").Value = 32
$ws2.Range("C28").Value = "ClkIn"
$ws2.Range("F28").Value = "DdsPdClk     "
$ws2.Range("G28").Value = "BS_Clk]      "
$ws2.Range("A29").Value = "BS_ADC_SDIO"
$ws2.Range("B29").Value = 34
$ws2.Range("C29").Value = "DdsData3 "
$ws2.Range("A30").Value = "BS_ADC_SClk"
$ws2.Range("B30").Value = 36
$ws2.Range("C30").Value = "DdsData2 "
$ws2.Range("A31").Value = "BS_ADC_CS_n"
$ws2.Range("B31").Value = 38
$ws2.Range("C31").Value = "DdsData1 "
$ws2.Range("A32").Value = "BS_ADC_LowZ"
$ws2.Range("B32").Value = 40
$ws2.Range("C32").Value = "DdsData0 "

# Column A best-fit width
$null = $ws2.Columns("A").AutoFit()

# Selection matching the saved view state
$null = $ws2.Range("A3:B32").Select()
